# Scheduled runner update: refresh Universalis market price / profit
# columns (H: currentAveragePrice, I: currentAveragePriceNQ,
# J: currentAveragePriceHQ, K: LevePriceNQ, L: LevePriceHQ,
# M: LeveProfitNQ, N: LeveProfitHQ) across the leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

# ALC row 28
$ALC.Range("H28").Value = 873.1429000000001
$ALC.Range("J28").Value = 1545
$ALC.Range("L28").Value = 1545
$ALC.Range("N28").Value = -2515

# ALC row 64
$ALC.Range("H64").Value = 13620
$ALC.Range("I64").Value = 9451
$ALC.Range("K64").Value = 9451
$ALC.Range("M64").Value = -9203

# ALC row 67
$ALC.Range("H67").Value = 13620
$ALC.Range("I67").Value = 9451
$ALC.Range("K67").Value = 9451
$ALC.Range("M67").Value = -8593

# ALC row 98
$ALC.Range("H98").Value = 1171.0344
$ALC.Range("I98").Value = 1044.4814
$ALC.Range("K98").Value = 1044.4814
$ALC.Range("M98").Value = 453.5186000000001

# ALC row 101
$ALC.Range("H101").Value = 1826.2
$ALC.Range("I101").Value = 1826.2
$ALC.Range("K101").Value = 5478.6
$ALC.Range("M101").Value = -3856.6

# ALC row 122
$ALC.Range("H122").Value = 1171.0344
$ALC.Range("I122").Value = 1044.4814
$ALC.Range("K122").Value = 3133.4442
$ALC.Range("M122").Value = -683.4441999999999

# ALC row 138
$ALC.Range("H138").Value = 10819.09
$ALC.Range("J138").Value = 10832.234
$ALC.Range("L138").Value = 32496.702
$ALC.Range("N138").Value = -42776.702

# ARM row 32
$ARM.Range("H32").Value = 4997197.5
$ARM.Range("I32").Value = 6612112
$ARM.Range("K32").Value = 6612112
$ARM.Range("M32").Value = -6611825

# ARM row 102
$ARM.Range("H102").Value = 58827708
$ARM.Range("I102").Value = 90914840
$ARM.Range("K102").Value = 90914840
$ARM.Range("M102").Value = -90913218

# ARM row 122
$ARM.Range("H122").Value = 2429.7058
$ARM.Range("I122").Value = 2429.7058
$ARM.Range("K122").Value = 7289.117400000001
$ARM.Range("M122").Value = -4839.117400000001

# CRP row 31
$CRP.Range("H31").Value = 4435.0557
$CRP.Range("I31").Value = 3741.4167
$CRP.Range("K31").Value = 3741.4167
$CRP.Range("M31").Value = -3446.4167

# CRP row 34
$CRP.Range("H34").Value = 4435.0557
$CRP.Range("I34").Value = 3741.4167
$CRP.Range("K34").Value = 3741.4167
$CRP.Range("M34").Value = -3539.4167

# CRP row 132
$CRP.Range("H132").Value = 9805827
$CRP.Range("I132").Value = 2107.4
$CRP.Range("K132").Value = 6322.200000000001
$CRP.Range("M132").Value = -3792.200000000001

# CRP row 134
$CRP.Range("H134").Value = 2746.0962
$CRP.Range("I134").Value = 1926.6216
$CRP.Range("K134").Value = 5779.864799999999
$CRP.Range("M134").Value = -3244.864799999999

# CRP row 141
$CRP.Range("H141").Value = 644030.8
$CRP.Range("J141").Value = 693036.25
$CRP.Range("L141").Value = 693036.25
$CRP.Range("N141").Value = -703396.25

# CUL row 2
$CUL.Range("H2").Value = 941.439
$CUL.Range("I2").Value = 312.3125
$CUL.Range("J2").Value = 1344.08
$CUL.Range("K2").Value = 1873.875
$CUL.Range("L2").Value = 8064.48
$CUL.Range("M2").Value = -1760.875
$CUL.Range("N2").Value = -8290.48

# CUL row 28
$CUL.Range("H28").Value = 9083.25
$CUL.Range("I28").Value = 2650
$CUL.Range("K28").Value = 7950
$CUL.Range("M28").Value = -7718

# CUL row 33
$CUL.Range("H33").Value = 127818.625
$CUL.Range("J33").Value = 204299.8
$CUL.Range("L33").Value = 1225798.8
$CUL.Range("N33").Value = -1226364.8

# CUL row 40
$CUL.Range("H40").Value = 2777.5
$CUL.Range("I40").Value = 5155
$CUL.Range("J40").Value = 400
$CUL.Range("K40").Value = 20620
$CUL.Range("L40").Value = 1600
$CUL.Range("M40").Value = -20551
$CUL.Range("N40").Value = -1738

# CUL row 88
$CUL.Range("H88").Value = 9291.242
$CUL.Range("I88").Value = 4778
$CUL.Range("J88").Value = 10097.179
$CUL.Range("K88").Value = 14334
$CUL.Range("L88").Value = 30291.537
$CUL.Range("M88").Value = -13906
$CUL.Range("N88").Value = -31147.537

# CUL row 91
$CUL.Range("H91").Value = 9291.242
$CUL.Range("I91").Value = 4778
$CUL.Range("J91").Value = 10097.179
$CUL.Range("K91").Value = 14334
$CUL.Range("L91").Value = 30291.537
$CUL.Range("M91").Value = -12852
$CUL.Range("N91").Value = -33255.537

# CUL row 97
$CUL.Range("H97").Value = 1494.6
$CUL.Range("J97").Value = 1927.3334
$CUL.Range("L97").Value = 5782.0002
$CUL.Range("N97").Value = -6774.0002

# CUL row 116
$CUL.Range("H116").Value = 3707.3845
$CUL.Range("I116").Value = 3719.6
$CUL.Range("K116").Value = 11158.8
$CUL.Range("M116").Value = -7716.799999999999

# CUL row 119
$CUL.Range("H119").Value = 125011360
$CUL.Range("I119").Value = 500006460
$CUL.Range("J119").Value = 12997.5
$CUL.Range("K119").Value = 1500019380
$CUL.Range("L119").Value = 38992.5
$CUL.Range("M119").Value = -1500014542
$CUL.Range("N119").Value = -48668.5

# CUL row 137
$CUL.Range("H137").Value = 11666.167
$CUL.Range("J137").Value = 11666.167
$CUL.Range("L137").Value = 34998.501
$CUL.Range("N137").Value = -45198.501

# GSM row 15
$GSM.Range("H15").Value = 54179.6
$GSM.Range("J15").Value = 61174.5
$GSM.Range("L15").Value = 61174.5
$GSM.Range("N15").Value = -61750.5

# GSM row 80
$GSM.Range("H80").Value = 33542228
$GSM.Range("I80").Value = 202925.62
$GSM.Range("J80").Value = 71644290
$GSM.Range("K80").Value = 202925.62
$GSM.Range("L80").Value = 71644290
$GSM.Range("M80").Value = -201927.62
$GSM.Range("N80").Value = -71646286

# GSM row 81
$GSM.Range("H81").Value = 54179.6
$GSM.Range("J81").Value = 61174.5
$GSM.Range("L81").Value = 61174.5
$GSM.Range("N81").Value = -63170.5

# GSM row 83
$GSM.Range("H83").Value = 33542228
$GSM.Range("I83").Value = 202925.62
$GSM.Range("J83").Value = 71644290
$GSM.Range("K83").Value = 1014628.1
$GSM.Range("L83").Value = 358221450
$GSM.Range("M83").Value = -1009636.1
$GSM.Range("N83").Value = -358231434

# GSM row 84
$GSM.Range("H84").Value = 54179.6
$GSM.Range("J84").Value = 61174.5
$GSM.Range("L84").Value = 183523.5
$GSM.Range("N84").Value = -193507.5

# GSM row 102
$GSM.Range("H102").Value = 71430000
$GSM.Range("I102").Value = 71430000
$GSM.Range("K102").Value = 71430000
$GSM.Range("M102").Value = -71428378

# GSM row 126
$GSM.Range("H126").Value = 5107.846
$GSM.Range("I126").Value = 4812.375
$GSM.Range("K126").Value = 14437.125
$GSM.Range("M126").Value = -11967.125

# GSM row 132
$GSM.Range("H132").Value = 6593425
$GSM.Range("I132").Value = 7590.5835
$GSM.Range("K132").Value = 22771.7505
$GSM.Range("M132").Value = -20241.7505

# LTW row 46
$LTW.Range("H46").Value = 10486.308
$LTW.Range("I46").Value = 11633.4
$LTW.Range("K46").Value = 11633.4
$LTW.Range("M46").Value = -11445.4

# LTW row 61
$LTW.Range("H61").Value = 1696.3903
$LTW.Range("I61").Value = 1576.3
$LTW.Range("K61").Value = 1576.3
$LTW.Range("M61").Value = -1374.3

# LTW row 113
$LTW.Range("H113").Value = 1696.3903
$LTW.Range("I113").Value = 1576.3
$LTW.Range("K113").Value = 1576.3
$LTW.Range("M113").Value = 593.7

# LTW row 139
$LTW.Range("H139").Value = 179863.17

# WVR row 52
$WVR.Range("H52").Value = 5000
$WVR.Range("I52").Value = 5000
$WVR.Range("K52").Value = 5000
$WVR.Range("M52").Value = -4774

# WVR row 107
$WVR.Range("H107").Value = 759.3226
$WVR.Range("I107").Value = 740.8570999999999
$WVR.Range("K107").Value = 2222.5713
$WVR.Range("M107").Value = -302.5712999999996

# WVR row 122
$WVR.Range("H122").Value = 10427
$WVR.Range("I122").Value = 10427
$WVR.Range("J122").Value = 0
$WVR.Range("K122").Value = 31281
$WVR.Range("L122").Value = 0
$WVR.Range("M122").Value = -28831
$WVR.Range("N122").ClearContents()  # cell removed: LevePriceHQ became 0, so no HQ profit value
